# Append the new "01-07-2021" monthly observation as row 32 on Sheet1,
# mirroring the existing rows 2-31 (Serie date in column A, then the
# Total/3 meses/6 meses/9 meses/12 meses/18 meses/2 años/5 años/10 años y
# más figures in columns B-J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 32

# Column A holds the period label as text (e.g. "01-06-2021"). Excel would
# otherwise auto-detect this "dd-mm-yyyy"-looking string and coerce it to a
# date serial, so force the cell to Text before typing it in, then drop the
# explicit format again so the cell keeps the sheet's default (unstyled)
# look, just like the rows above it.
$ws.Range("A" + $row).NumberFormat = "@"
$ws.Range("A" + $row).Value = "01-07-2021"
$ws.Range("A" + $row).Style = "Normal"

$values = @{
    "B" = 78499
    "C" = 12881
    "D" = 9312
    "E" = 4983
    "F" = 5181
    "G" = 6269
    "H" = 17107
    "I" = 13890
    "J" = 8875
}

foreach ($col in $values.Keys) {
    $ws.Range($col + $row).Value = $values[$col]
}
